# "Added example w.r.t 20.10"
#
# This commit does two things to the single-slide example deck:
#   1. Resizes the picture on slide 1 (it had been inserted at a
#      placeholder 6096000 x 6096000 EMU square; it's resized to its
#      actual/intended extent).
#   2. Adds a review comment thread (a top-level comment plus a reply)
#      anchored near the top-left corner of the slide.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- 1. Resize the picture -------------------------------------------------
$pic = $s.Shapes.Item(1)
# 3707904 EMU / 12700 EMU-per-point, 1124744 EMU / 12700 EMU-per-point
$pic.Width = 291.9609680175781
$pic.Height = 88.56251968503938

# --- 2. Add the comment thread ---------------------------------------------
# Comments.Add(Left, Top, Author, AuthorInitials, Text) takes Left/Top in
# points; the target position is (10, 10) EMU, i.e. a tiny fraction of a
# point from the top-left corner.
$comment = $s.Comments.Add(0.0007874015748031496, 0.0007874015748031496, "Сергей Пучок", "СП", "The comment")

# Reply, threaded under the comment above.
$reply = $comment.Replies.Add("Сергей Пучок", "СП", "Reply comment")
